# Generate Report for Handback
# Update handoff/handback timestamps for the first file row on the
# zh-cn and de-de localization status sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-01 22:53:10"
$wsZhCn.Range("K2").Value = "2016-09-01 22:53:28"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-01 22:53:15"
$wsDeDe.Range("K2").Value = "2016-09-01 22:53:35"
